$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's row (2025-05-21) at the bottom of the "Prices" table.
# Carries forward the previous day's Argent/Silver price values unchanged,
# matching the source data which stores every cell as literal text.
$row = 81
$ws.Range("A$row").Value = "'2025-05-21"
$ws.Range("B$row").Value = "'36.5"
$ws.Range("C$row").Value = "'36.25"
$ws.Range("D$row").Value = "'0.94"
$ws.Range("E$row").Value = "'0.258"
$ws.Range("F$row").Value = "'0.09"
$ws.Range("G$row").Value = "'5,289"
$ws.Range("H$row").Value = "'7,918"
$ws.Range("I$row").Value = "'7,968"
$ws.Range("J$row").Value = "'7.2326"
